$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 9
$ws.Range("A3").Value = 9
$ws.Range("A4").Value = 9

$ws.Range("A5").Value = 14
$ws.Range("A6").Value = 14
$ws.Range("A7").Value = 14

$ws.Range("A8").Value = 20
$ws.Range("A9").Value = 20
$ws.Range("A10").Value = 20
